$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '97.020.41'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.710.35'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.38%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.82%  '
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '656.23'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.88%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.430'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.88%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.07'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.63%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.707.93'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '44.37'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.71%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.208'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.08%  '
$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000304'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +13.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.78'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.401.20'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '96.708.21'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.92'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.697.87'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.18'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '18.73'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.508'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '523.20'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.42'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.32%  '
$ws.Range('E25').Value = '  +1.43%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.95'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.25%  '
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.195'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +16.47%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '13.38'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.43%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '12.19'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.21%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.01'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.81%  '
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('E33').Value = '  +1.19%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.88'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +7.86%  '
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '32.36'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '652.80'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.46%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.599'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.50%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.85'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.46%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.82'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +10.38%  '
$ws.Range('B42').Value = 'ImmutableX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.04'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.01%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.84'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.161'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.94%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.964'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.18%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.446'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.60%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0455'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.57%  '
$ws.Range('E48').Value = '  -0.60%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '23.62'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.54'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.68%  '
$ws.Range('E51').Value = '  +0.47%  '
